# Weekly crime data refresh: update report dates/volume number and all
# period-over-period statistics for the new reporting week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
# A8 = "Volume 32   Number  30" -> "...31"
$a8 = $ws.Cells.Item(8, 1)
$a8.Characters(21, 2).Text = "31"

# C9 = "Report Covering the Week  7/21/2025  Through  7/27/2025"
#   -> "Report Covering the Week  7/28/2025  Through  8/3/2025"
$c9 = $ws.Cells.Item(9, 3)
$c9.Characters(47, 9).Text = "8/3/2025"
$c9.Characters(27, 9).Text = "7/28/2025"

# --- Statistical table updates (rows 15-31) ---
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 8).Value = 100
$ws.Cells.Item(15, 13).Value = 420
$ws.Cells.Item(15, 14).Value = 8.333333333333
$ws.Cells.Item(16, 3).Value = 4
$ws.Cells.Item(16, 4).Value = 11
$ws.Cells.Item(16, 5).Value = -63.636363636363
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(16, 7).Value = 32
$ws.Cells.Item(16, 8).Value = -53.125
$ws.Cells.Item(16, 9).Value = 193
$ws.Cells.Item(16, 10).Value = 236
$ws.Cells.Item(16, 11).Value = -18.220338983050
$ws.Cells.Item(16, 12).Value = -31.071428571428
$ws.Cells.Item(16, 13).Value = 129.761904761905
$ws.Cells.Item(16, 14).Value = -87.302631578947
$ws.Cells.Item(17, 3).Value = 13
$ws.Cells.Item(17, 4).Value = 18
$ws.Cells.Item(17, 5).Value = -27.777777777777
$ws.Cells.Item(17, 6).Value = 45
$ws.Cells.Item(17, 7).Value = 51
$ws.Cells.Item(17, 8).Value = -11.764705882352
$ws.Cells.Item(17, 9).Value = 337
$ws.Cells.Item(17, 10).Value = 330
$ws.Cells.Item(17, 11).Value = 2.121212121212
$ws.Cells.Item(17, 12).Value = 13.087248322147
$ws.Cells.Item(17, 13).Value = 195.614035087719
$ws.Cells.Item(17, 14).Value = -9.164420485175
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(18, 5).Value = 233.333333333333
$ws.Cells.Item(18, 6).Value = 40
$ws.Cells.Item(18, 8).Value = 60
$ws.Cells.Item(18, 9).Value = 252
$ws.Cells.Item(18, 10).Value = 218
$ws.Cells.Item(18, 11).Value = 15.596330275229
$ws.Cells.Item(18, 12).Value = -3.076923076923
$ws.Cells.Item(18, 13).Value = 28.571428571428
$ws.Cells.Item(18, 14).Value = -83.752417794971
$ws.Cells.Item(19, 3).Value = 30
$ws.Cells.Item(19, 4).Value = 29
$ws.Cells.Item(19, 5).Value = 3.448275862068
$ws.Cells.Item(19, 6).Value = 112
$ws.Cells.Item(19, 7).Value = 161
$ws.Cells.Item(19, 8).Value = -30.434782608695
$ws.Cells.Item(19, 9).Value = 1021
$ws.Cells.Item(19, 10).Value = 1209
$ws.Cells.Item(19, 11).Value = -15.550041356493
$ws.Cells.Item(19, 12).Value = -27.588652482269
$ws.Cells.Item(19, 13).Value = -24.370370370370
$ws.Cells.Item(19, 14).Value = -81.603603603603
$ws.Cells.Item(20, 7).Value = 3
$ws.Cells.Item(20, 8).Value = -66.666666666666
$ws.Cells.Item(20, 10).Value = 31
$ws.Cells.Item(20, 11).Value = -58.064516129032
$ws.Cells.Item(20, 12).Value = -65.789473684210
$ws.Cells.Item(20, 14).Value = -93.779904306220
$ws.Cells.Item(21, 3).Value = 57
$ws.Cells.Item(21, 4).Value = 62
$ws.Cells.Item(21, 5).Value = -8.064516129032
$ws.Cells.Item(21, 6).Value = 215
$ws.Cells.Item(21, 7).Value = 274
$ws.Cells.Item(21, 8).Value = -21.532846715328
$ws.Cells.Item(21, 9).Value = 1843
$ws.Cells.Item(21, 10).Value = 2047
$ws.Cells.Item(21, 11).Value = -9.965803615046
$ws.Cells.Item(21, 12).Value = -19.764910753156
$ws.Cells.Item(21, 13).Value = 4.537719795802
$ws.Cells.Item(21, 14).Value = -80.036828422877

# C22 changes from a blank-marker ("0") text cell to a real numeric cell;
# match the numeric format used by its neighboring data cells.
$ws.Cells.Item(22, 3).Value = 3
$ws.Cells.Item(22, 3).NumberFormat = "#,##0"
$ws.Cells.Item(22, 4).Value = 4
$ws.Cells.Item(22, 5).Value = -25
$ws.Cells.Item(22, 6).Value = 9
$ws.Cells.Item(22, 7).Value = 12
$ws.Cells.Item(22, 8).Value = -25
$ws.Cells.Item(22, 9).Value = 123
$ws.Cells.Item(22, 10).Value = 106
$ws.Cells.Item(22, 11).Value = 16.037735849056
$ws.Cells.Item(22, 12).Value = -6.106870229007
$ws.Cells.Item(22, 13).Value = 48.192771084337
$ws.Cells.Item(24, 3).Value = 84
$ws.Cells.Item(24, 4).Value = 100
$ws.Cells.Item(24, 5).Value = -16
$ws.Cells.Item(24, 6).Value = 323
$ws.Cells.Item(24, 7).Value = 363
$ws.Cells.Item(24, 8).Value = -11.019283746556
$ws.Cells.Item(24, 9).Value = 2381
$ws.Cells.Item(24, 10).Value = 2588
$ws.Cells.Item(24, 11).Value = -7.998454404945
$ws.Cells.Item(24, 12).Value = 0.168279343710
$ws.Cells.Item(24, 13).Value = -12.623853211009
$ws.Cells.Item(25, 3).Value = 75
$ws.Cells.Item(25, 4).Value = 88
$ws.Cells.Item(25, 5).Value = -14.772727272727
$ws.Cells.Item(25, 6).Value = 272
$ws.Cells.Item(25, 7).Value = 320
$ws.Cells.Item(25, 8).Value = -15
$ws.Cells.Item(25, 9).Value = 2040
$ws.Cells.Item(25, 10).Value = 2290
$ws.Cells.Item(25, 11).Value = -10.917030567685
$ws.Cells.Item(25, 12).Value = -7.019143117593
$ws.Cells.Item(26, 3).Value = 16
$ws.Cells.Item(26, 4).Value = 31
$ws.Cells.Item(26, 5).Value = -48.387096774193
$ws.Cells.Item(26, 6).Value = 93
$ws.Cells.Item(26, 7).Value = 87
$ws.Cells.Item(26, 8).Value = 6.896551724137
$ws.Cells.Item(26, 9).Value = 652
$ws.Cells.Item(26, 10).Value = 616
$ws.Cells.Item(26, 11).Value = 5.844155844155
$ws.Cells.Item(26, 12).Value = 2.677165354330
$ws.Cells.Item(26, 13).Value = 95.209580838323
$ws.Cells.Item(27, 6).Value = 3
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 200
$ws.Cells.Item(27, 12).Value = 70.588235294117
$ws.Cells.Item(28, 3).Value = 15
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).Value = 650
$ws.Cells.Item(28, 6).Value = 58
$ws.Cells.Item(28, 7).Value = 21
$ws.Cells.Item(28, 8).Value = 176.190476190476
$ws.Cells.Item(28, 9).Value = 197
$ws.Cells.Item(28, 10).Value = 132
$ws.Cells.Item(28, 11).Value = 49.242424242424
$ws.Cells.Item(28, 12).Value = 45.925925925925
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(31, 7).Value = 6
$ws.Cells.Item(31, 8).Value = -83.333333333333
$ws.Cells.Item(31, 9).Value = 11
$ws.Cells.Item(31, 10).Value = 18
$ws.Cells.Item(31, 11).Value = -38.888888888888
$ws.Cells.Item(31, 12).Value = 10
